$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Cells.Item(4, 2).Value = "inf"
$ws.Cells.Item(6, 2).Value = 733671.9799030328
$ws.Cells.Item(7, 2).Value = 1784511.217947469
$ws.Cells.Item(8, 2).Value = 19042283.03931422
$ws.Cells.Item(10, 2).Value = 6746623.038510369

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Cells.Item(2, 2).Value = 546846.2755035073
$ws.Cells.Item(2, 3).Value = 546846.5903286961
$ws.Cells.Item(2, 4).Value = 546853.5913296627
$ws.Cells.Item(2, 5).Value = 183139.6229770186
$ws.Cells.Item(2, 6).Value = 183139.6229770186
$ws.Cells.Item(2, 7).Value = 183139.6229770186
$ws.Cells.Item(2, 8).Value = 183139.6229770186
$ws.Cells.Item(2, 9).Value = 183139.6229770186
$ws.Cells.Item(2, 10).Value = 183139.6229770186
$ws.Cells.Item(2, 11).Value = 183139.6229770186
$ws.Cells.Item(2, 12).Value = 183139.6229770186
$ws.Cells.Item(2, 13).Value = 183139.6229770186
$ws.Cells.Item(2, 14).Value = 183139.6229770186
$ws.Cells.Item(2, 15).Value = 183139.6229770186
$ws.Cells.Item(2, 16).Value = 183139.6229770186
$ws.Cells.Item(3, 2).Value = 329223.0061167778
$ws.Cells.Item(3, 3).Value = 1558.045904853998
$ws.Cells.Item(3, 4).Value = 32585.62631485736
$ws.Cells.Item(4, 2).Value = 418430.1700736278
$ws.Cells.Item(4, 3).Value = 417860.3230835526
$ws.Cells.Item(4, 5).Value = 19822.11383057336
$ws.Cells.Item(4, 6).Value = 19822.11383057336
$ws.Cells.Item(4, 7).Value = 19822.11383057336
$ws.Cells.Item(4, 9).Value = 19822.11383057336
$ws.Cells.Item(4, 10).Value = 19822.11383057336
$ws.Cells.Item(4, 11).Value = 19822.11383057336
$ws.Cells.Item(4, 12).Value = 19822.11383057336
$ws.Cells.Item(4, 13).Value = 19822.11383057336
$ws.Cells.Item(4, 14).Value = 19822.11383057336
$ws.Cells.Item(4, 15).Value = 19822.11383057336
$ws.Cells.Item(4, 16).Value = 19822.11383057336
$ws.Cells.Item(5, 2).Value = 41075.67462471527
$ws.Cells.Item(6, 2).Value = -241882.5753116135
$ws.Cells.Item(6, 3).Value = 86315.2467747375
$ws.Cells.Item(6, 4).Value = 66536.26763429289
$ws.Cells.Item(6, 5).Value = 83768.00283179304
$ws.Cells.Item(6, 6).Value = 153120.8091476814
$ws.Cells.Item(6, 7).Value = 153120.8091476814
$ws.Cells.Item(6, 8).Value = 153120.8091476814
$ws.Cells.Item(6, 9).Value = 153120.8091476814
$ws.Cells.Item(6, 10).Value = 153120.8091476814
$ws.Cells.Item(6, 11).Value = 153120.8091476814
$ws.Cells.Item(6, 12).Value = 153120.8091476814
$ws.Cells.Item(6, 13).Value = 153120.8091476814
$ws.Cells.Item(6, 14).Value = 153120.8091476814
$ws.Cells.Item(6, 15).Value = 153120.8091476814
$ws.Cells.Item(6, 16).Value = 153120.8091476814

$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Cells.Item(3, 2).Value = 341.6547992988656

$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Cells.Item(3, 2).Value = 341.6547992988656
$ws.Cells.Item(3, 3).Value = 1.71100646040092
$ws.Cells.Item(3, 4).Value = 38.04891829752592

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Cells.Item(2, 8).Value = 350.0302086235789
$ws.Cells.Item(2, 9).Value = 254.6788077461599
$ws.Cells.Item(2, 10).Value = 104.241554533109
$ws.Cells.Item(2, 11).Value = 62.89094992975379
$ws.Cells.Item(2, 12).Value = 19.46997469197376
$ws.Cells.Item(2, 15).Value = 6.093541586991591
$ws.Cells.Item(2, 16).Value = 57.00149786528118
$ws.Cells.Item(2, 17).Value = 121.5990695169012
$ws.Cells.Item(2, 18).Value = 211.9303003426456
$ws.Cells.Item(2, 19).Value = 236.0914234572714
$ws.Cells.Item(2, 20).Value = 219.5019606424398
$ws.Cells.Item(2, 21).Value = 248.7917266063122
$ws.Cells.Item(3, 7).Value = 161.0343157549836
$ws.Cells.Item(3, 8).Value = 137.937230470563
$ws.Cells.Item(3, 9).Value = 117.6169923923482
$ws.Cells.Item(3, 10).Value = 83.28630290869316
$ws.Cells.Item(3, 11).Value = 14.16511697308491
$ws.Cells.Item(3, 17).Value = 56.29644527255945
$ws.Cells.Item(3, 18).Value = 155.5936088975356
$ws.Cells.Item(3, 19).Value = 208.3819030908814
$ws.Cells.Item(3, 20).Value = 230.1732865992993
$ws.Cells.Item(3, 21).Value = 249.6508274585269
$ws.Cells.Item(4, 7).Value = 169.250603464539
$ws.Cells.Item(4, 8).Value = 167.7218531446305
$ws.Cells.Item(4, 9).Value = 168.2233735905427
$ws.Cells.Item(4, 10).Value = 133.427194649778
$ws.Cells.Item(4, 11).Value = 89.34902992251145
$ws.Cells.Item(4, 12).Value = 63.42832505831589
$ws.Cells.Item(4, 13).Value = 60.19015695232082
$ws.Cells.Item(4, 14).Value = 47.78758048940401
$ws.Cells.Item(4, 15).Value = 73.019086697924
$ws.Cells.Item(4, 16).Value = 91.81114782298006
$ws.Cells.Item(4, 17).Value = 152.9970210107205
$ws.Cells.Item(4, 18).Value = 221.7828567288786
$ws.Cells.Item(4, 19).Value = 243.9747107992991
$ws.Cells.Item(4, 20).Value = 218.1559507844908
$ws.Cells.Item(4, 21).Value = 291.2204648267885
$ws.Cells.Item(5, 8).Value = 349.9597650259403
$ws.Cells.Item(5, 9).Value = 254.413627538865
$ws.Cells.Item(5, 10).Value = 103.6577574092157
$ws.Cells.Item(5, 11).Value = 62.01598937483024
$ws.Cells.Item(5, 12).Value = 18.38450875428637
$ws.Cells.Item(5, 15).Value = 4.934605547775959
$ws.Cells.Item(5, 16).Value = 56.01237276874352
$ws.Cells.Item(5, 17).Value = 120.8562777625268
$ws.Cells.Item(5, 18).Value = 211.4982239222711
$ws.Cells.Item(5, 19).Value = 235.9346815086628
$ws.Cells.Item(5, 20).Value = 219.4718503679457
$ws.Cells.Item(5, 21).Value = 248.7911763328777
$ws.Cells.Item(6, 7).Value = 161.0306354769367
$ws.Cells.Item(6, 8).Value = 137.9016867325838
$ws.Cells.Item(6, 9).Value = 117.4902810648562
$ws.Cells.Item(6, 10).Value = 82.93859734111298
$ws.Cells.Item(6, 11).Value = 13.57083277636241
$ws.Cells.Item(6, 17).Value = 55.82666100818672
$ws.Cells.Item(6, 18).Value = 155.3651088272201
$ws.Cells.Item(6, 19).Value = 208.3135435403172
$ws.Cells.Item(6, 20).Value = 230.1584524961191
$ws.Cells.Item(6, 21).Value = 249.6505853349712
$ws.Cells.Item(7, 7).Value = 169.247518043053
$ws.Cells.Item(7, 8).Value = 167.6944209426916
$ws.Cells.Item(7, 9).Value = 168.1305865516738
$ws.Cells.Item(7, 10).Value = 133.20905535072
$ws.Cells.Item(7, 11).Value = 88.99056004441432
$ws.Cells.Item(7, 12).Value = 62.96960703121101
$ws.Cells.Item(7, 13).Value = 59.70650310975206
$ws.Cells.Item(7, 14).Value = 47.31542685419204
$ws.Cells.Item(7, 15).Value = 72.58297639552539
$ws.Cells.Item(7, 16).Value = 91.43798011889523
$ws.Cells.Item(7, 17).Value = 152.7386590351999
$ws.Cells.Item(7, 18).Value = 221.644124959155
$ws.Cells.Item(7, 19).Value = 243.9209403175845
$ws.Cells.Item(7, 20).Value = 218.1427676199599
$ws.Cells.Item(7, 21).Value = 291.2202965310711

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Cells.Item(2, 7).Value = 1.373486630347197
$ws.Cells.Item(2, 8).Value = 14.06621995304324
$ws.Cells.Item(2, 9).Value = 52.95134331646037
$ws.Cells.Item(2, 10).Value = 116.5729608924305
$ws.Cells.Item(2, 11).Value = 174.7126499550274
$ws.Cells.Item(2, 12).Value = 216.7464914185155
$ws.Cells.Item(2, 13).Value = 241.1722342809525
$ws.Cells.Item(2, 14).Value = 245.0746531694265
$ws.Cells.Item(2, 15).Value = 231.4170454889115
$ws.Cells.Item(2, 16).Value = 197.509094302215
$ws.Cells.Item(2, 17).Value = 148.321104352906
$ws.Cells.Item(2, 18).Value = 86.27727954354719
$ws.Cells.Item(2, 19).Value = 31.29832658903679
$ws.Cells.Item(2, 20).Value = 6.012437724344859
$ws.Cells.Item(2, 21).Value = 0.1098789304277757
$ws.Cells.Item(3, 7).Value = 0.7348801343409561
$ws.Cells.Item(3, 8).Value = 7.09739498166134
$ws.Cells.Item(3, 9).Value = 25.30179409901977
$ws.Cells.Item(3, 10).Value = 69.4300569028006
$ws.Cells.Item(3, 11).Value = 118.6670259036447
$ws.Cells.Item(3, 12).Value = 159.5624607480221
$ws.Cells.Item(3, 13).Value = 186.2018656178817
$ws.Cells.Item(3, 14).Value = 191.1300749398437
$ws.Cells.Item(3, 15).Value = 174.8466782789204
$ws.Cells.Item(3, 16).Value = 140.3298740742831
$ws.Cells.Item(3, 17).Value = 93.80680451692626
$ws.Cells.Item(3, 18).Value = 45.62703149881833
$ws.Cells.Item(3, 19).Value = 13.65007617953485
$ws.Cells.Item(3, 20).Value = 2.962082646751484
$ws.Cells.Item(3, 21).Value = 0.04834737725927345
$ws.Cells.Item(4, 7).Value = 0.6160988184077902
$ws.Cells.Item(4, 8).Value = 5.477678585480175
$ws.Cells.Item(4, 9).Value = 18.52777173902701
$ws.Cells.Item(4, 10).Value = 43.55818646143077
$ws.Cells.Item(4, 11).Value = 71.57948090228689
$ws.Cells.Item(4, 12).Value = 91.59709160219094
$ws.Cells.Item(4, 13).Value = 96.57629023459569
$ws.Cells.Item(4, 14).Value = 94.27992191143946
$ws.Cells.Item(4, 15).Value = 87.08276753276661
$ws.Cells.Item(4, 16).Value = 74.51435163724763
$ws.Cells.Item(4, 17).Value = 51.5898746941287
$ws.Cells.Item(4, 18).Value = 27.70204323495391
$ws.Cells.Item(4, 19).Value = 10.73692213534303
$ws.Cells.Item(4, 20).Value = 2.632422224106012
$ws.Cells.Item(4, 21).Value = 0.03360539009497042
$ws.Cells.Item(5, 7).Value = 1.380365048278458
$ws.Cells.Item(5, 8).Value = 14.13666355068176
$ws.Cells.Item(5, 9).Value = 53.21652352375528
$ws.Cells.Item(5, 10).Value = 117.1567580163238
$ws.Cells.Item(5, 11).Value = 175.587610509951
$ws.Cells.Item(5, 12).Value = 217.8319573562029
$ws.Cells.Item(5, 13).Value = 242.3800242835249
$ws.Cells.Item(5, 14).Value = 246.301986476946
$ws.Cells.Item(5, 15).Value = 232.5759815281271
$ws.Cells.Item(5, 16).Value = 198.4982193987527
$ws.Cells.Item(5, 17).Value = 149.0638961072804
$ws.Cells.Item(5, 18).Value = 86.70935596392172
$ws.Cells.Item(5, 19).Value = 31.45506853764538
$ws.Cells.Item(5, 20).Value = 6.04254799883895
$ws.Cells.Item(5, 21).Value = 0.1104292038622766
$ws.Cells.Item(6, 7).Value = 0.7385604123878564
$ws.Cells.Item(6, 8).Value = 7.132938719640613
$ws.Cells.Item(6, 9).Value = 25.42850542651172
$ws.Cells.Item(6, 10).Value = 69.77776247038078
$ws.Cells.Item(6, 11).Value = 119.2613101003672
$ws.Cells.Item(6, 12).Value = 160.3615491897405
$ws.Cells.Item(6, 13).Value = 187.1343641388002
$ws.Cells.Item(6, 15).Value = 175.7223100473861
$ws.Cells.Item(6, 16).Value = 141.0326457655372
$ws.Cells.Item(6, 17).Value = 94.276588781299
$ws.Cells.Item(6, 18).Value = 45.85553156913376
$ws.Cells.Item(6, 19).Value = 13.71843573009899
$ws.Cells.Item(6, 20).Value = 2.976916749931753
$ws.Cells.Item(6, 21).Value = 0.04858950081499057
$ws.Cells.Item(7, 7).Value = 0.6191842398937593
$ws.Cells.Item(7, 8).Value = 5.505110787419063
$ws.Cells.Item(7, 9).Value = 18.62055877789597
$ws.Cells.Item(7, 10).Value = 43.77632576048878
$ws.Cells.Item(7, 11).Value = 71.93795078038401
$ws.Cells.Item(7, 12).Value = 92.05580962929582
$ws.Cells.Item(7, 13).Value = 97.05994407716445
$ws.Cells.Item(7, 14).Value = 94.75207554665143
$ws.Cells.Item(7, 15).Value = 87.51887783516521
$ws.Cells.Item(7, 16).Value = 74.88751934133246
$ws.Cells.Item(7, 17).Value = 51.84823666964925
$ws.Cells.Item(7, 18).Value = 27.84077500467757
$ws.Cells.Item(7, 19).Value = 10.7906926170576
$ws.Cells.Item(7, 20).Value = 2.645605388636971
$ws.Cells.Item(7, 21).Value = 0.03377368581238691
$ws.Cells.Item(20, 17).Value = 203.0572840332874
$ws.Cells.Item(20, 20).Value = 8.231257986185739
$ws.Cells.Item(21, 21).Value = 0.06618941491321523
$ws.Cells.Item(22, 8).Value = 7.499152202429723
$ws.Cells.Item(22, 11).Value = 97.9950563875785

$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Cells.Item(2, 13).Value = 95.70913181393388
$ws.Cells.Item(2, 14).Value = 90.44195234058779
$ws.Cells.Item(3, 12).Value = 43.07331968004391
$ws.Cells.Item(3, 13).Value = 21.66915899307011
$ws.Cells.Item(3, 14).Value = 0.9571789820313086
$ws.Cells.Item(3, 15).Value = 33.70032922107961
$ws.Cells.Item(3, 16).Value = 55.60769676917488
$ws.Cells.Item(5, 13).Value = 94.50134181136147
$ws.Cells.Item(5, 14).Value = 89.2146190330682
$ws.Cells.Item(5, 15).Value = 99.00804712831379
$ws.Cells.Item(6, 12).Value = 42.27423123832551
$ws.Cells.Item(6, 13).Value = 20.73666047215158
$ws.Cells.Item(6, 15).Value = 32.82469745261383
$ws.Cells.Item(6, 16).Value = 54.90492507792075
$ws.Cells.Item(20, 11).Value = 82.693084352536

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Cells.Item(2, 7).Value = 22.21673494391485
$ws.Cells.Item(5, 7).Value = 22.20985652598358
$ws.Cells.Item(8, 7).Value = 22.05689605041567
$ws.Cells.Item(20, 17).Value = 66.8628898365198
$ws.Cells.Item(22, 11).Value = 62.93345443721984
$ws.Cells.Item(22, 12).Value = 29.62548300147245

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Cells.Item(2, 2).Value = 1154526.17053041
$ws.Cells.Item(3, 2).Value = 1154339.431969725
